# Reorders the player roster rows (A2:C19) back to the "draft order" layout.
# Player identity (column A) moves together with its Position (B) and Team (C)
# as a single unit - only the row order changes, no values are altered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dyson Daniels",        "PG,SG",    "Atlanta Hawks"),
    @("Tari Eason",           "SF,PF",    "Houston Rockets"),
    @("De'Andre Hunter",      "SF,PF",    "Atlanta Hawks"),
    @("Jamal Murray",         "PG,SG",    "Denver Nuggets"),
    @("Victor Wembanyama",    "C",        "San Antonio Spurs"),
    @("Myles Turner",         "C",        "Indiana Pacers"),
    @("Deandre Ayton",        "C",        "Portland Trail Blazers"),
    @("Bradley Beal",         "PG,SG,SF", "Phoenix Suns"),
    @("Donovan Mitchell",     "PG,SG",    "Cleveland Cavaliers"),
    @("Michael Porter Jr.",   "SF,PF",    "Denver Nuggets"),
    @("Malik Beasley",        "SG",       "Detroit Pistons"),
    @("Payton Pritchard",     "PG",       "Boston Celtics"),
    @("Kristaps Porzingis",   "PF,C",     "Boston Celtics"),
    @("Domantas Sabonis",     "C",        "Sacramento Kings"),
    @("Josh Hart",            "SF,PF",    "New York Knicks"),
    @("Robert Williams III",  "C",        "Portland Trail Blazers"),
    @("Cam Thomas",           "SG,SF",    "Brooklyn Nets"),
    @("P.J. Washington",      "PF",       "Dallas Mavericks")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
